$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (A1, C1); B1 "Location" stays the same.
$ws.Range("A1").Value = "ParkName"
$ws.Range("C1").Value = "NumberOfVisitors"

# Set widths for the two newly introduced (still empty) columns D and E.
$ws.Columns("D:D").ColumnWidth = 14.6
$ws.Columns("E:E").ColumnWidth = 15.25

# Update the saved selection/active cell to E13.
$ws.Range("E13").Select()
